# bug fix in Eduati data files
#
# Sheet1 ("COLO320HSR_noCTRL_meas.xlsx") had 43 stray trailing rows
# (rows 45:87) that only held a leftover index column (A) with no real
# measurement data - remove them, which also shrinks the used range back
# down to A1:N44. The previously-active sheet (Sheet3) loses focus/selection
# in favour of Sheet1, and view selections are refreshed to reflect where
# each sheet was last clicked.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet3: was the active/selected tab before the fix; record its new
#     (non-active) selection before we move focus away from it.
$ws3.Activate()
[void]$ws3.Range("B43").Select()

# --- Sheet1: drop the bogus extra rows 45:87 (index-only leftovers beyond
#     the real 44 rows of data), which pulls the sheet's dimension back to
#     A1:N44.
[void]$ws1.Range("A45:A87").EntireRow.Delete()

# --- Sheet1 becomes the active sheet/tab again, scrolled down towards the
#     bottom of the (now shorter) data and with cell F62 selected.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("F62").Select()

# Sheet2 is untouched by this fix.
